$wb = $excel.ActiveWorkbook

# "Poland" is the last sheet/tab in the gallery; the new "UK" market sheet
# is added the same way a user would in Excel: duplicate the most similar
# existing sheet (Poland) and edit the two market-specific cells.
$poland = $wb.Worksheets.Item("Poland")
$poland.Copy($null, $poland)

# The copy lands immediately after "Poland" and becomes the active sheet.
$uk = $wb.Worksheets.Item($wb.Worksheets.Count)
$uk.Name = "UK"

# Fill in the market-specific values (User Story first, then Description,
# matching the authoring order reflected in the shared-string table).
$uk.Range("B4").Value = "NGC-2741/T3334"
$uk.Range("B2").Value = "UK Market"

# Leave the new sheet's selection on the User Story cell, and make it the
# active/visible tab.
$uk.Range("B4").Select()
$uk.Activate()
